$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(33, 1).Value = 10002
$ws.Cells.Item(33, 2).Value = 110032
$ws.Cells.Item(33, 3).Value = 10032
$ws.Cells.Item(33, 4).Value = "eng"
$ws.Cells.Item(33, 5).Value = $true
$ws.Cells.Item(33, 6).Value = "superadmin"
$ws.Cells.Item(33, 7).Value = "now()"

$ws.Range("C30").Select()
